# CCO_eCoaching_Log_Runbook.docx edit script
# Implements: "#18448 - CSR Incentive Data Feed" commit
#   - TFS 18175 / "Allow senior managers..." entry -> TFS 18448 / "CSR Incentive Data Feed (IDD)"
#   - New log row: 09/24/2020 | TFS 18448 - CSR Incentive Data Feed (IDD) | Lili Huang
#   - Purpose statement TFS number 18321 -> 18448
#   - Changeset number 46949 -> 47208
#   - web.config step: restructure the Prod.VnV.IPs instructions
#
$d = $word.ActiveDocument

# Useful characters for curly punctuation used throughout this document.
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”
$rsq = [char]0x2019   # '
$ndash = [char]0x2013 # -

# ---------------------------------------------------------------
# 1) "TFS 18175 - Allow senior managers to view log details."
#    becomes "TFS 18448 - CSR Incentive Data Feed (IDD)"
# ---------------------------------------------------------------
$old1 = "TFS 18175 " + $ndash + " Allow senior managers to view log details."
$new1 = "TFS 18448 " + $ndash + " CSR Incentive Data Feed (IDD)"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) { throw "Step 1: could not find the TFS 18175 summary text" }

# ---------------------------------------------------------------
# 2) Add a new row to the change-log table (Date / Change Description / Author)
#    right after the existing "08/28/2020 - TFS 18321 ..." row.
# ---------------------------------------------------------------
$logTableIndex = -1
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Columns.Count -eq 3 -and $candidate.Rows.Count -gt 1) {
        $headerText = $candidate.Rows.Item(1).Range.Text
        if ($headerText -like "Date*Change Description*Author*") {
            $logTableIndex = $i
        }
    }
}

if ($logTableIndex -eq -1) { throw "Could not locate the Date/Change Description/Author log table" }

$logTable = $d.Tables.Item($logTableIndex)
$newRow = $logTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "09/24/2020"
$newRow.Cells.Item(2).Range.Text = "TFS 18448 " + $ndash + " CSR Incentive Data Feed (IDD)"
$newRow.Cells.Item(3).Range.Text = "Lili Huang"

# ---------------------------------------------------------------
# 3) Purpose statement: "... per TFS 18321." -> "... per TFS 18448."
# ---------------------------------------------------------------
$old3 = "per TFS 18321."
$new3 = "per TFS 18448."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $found3) { throw "Step 3: could not find the purpose statement TFS number" }

# ---------------------------------------------------------------
# 4) Changeset number: 46949 -> 47208
# ---------------------------------------------------------------
$old4 = "Changeset 46949"
$new4 = "Changeset 47208"
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $found4) { throw "Step 4: could not find the Changeset number" }

# ---------------------------------------------------------------
# 5) web.config step: rework the "Prod.VnV.IPs" instruction sentence,
#    keeping "Prod.VnV.IPs" itself bold/untouched.
#    Before: "<ldq>Prod.VnV.IPs<rdq>, add/update the testers<rsq> IPs, separated by comma;"
#    After:  "<ldq>Prod.VnV.IPs<rdq>,  remove all exiting IPs, then add testers<rsq> IPs, separated by comma;"
# ---------------------------------------------------------------
$old5 = $rdq + ", add/update the testers" + $rsq + " I"
$new5 = $rdq + ",  remove all exiting IPs, then add testers" + $rsq + " I"
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $found5) { throw "Step 5: could not find the Prod.VnV.IPs instruction text" }

Write-Output "done"
